$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('A big plus EOXS works out of the box with mill level specs test certs and even processing add ons', 0, 1, 1, 100),
    @('Appreciate the help EOXS might just be our match', 0, 1, 1, 100),
    @('Can EOXS handle processing operations too', 0, 1, 1, 100),
    @('Can you explain the key features of an ERP system designed for metal distribution companies', 6, 0, 6, 0),
    @('Do you have more details on EOXS s capabilities', 0, 1, 1, 100),
    @('EOXS I ve seen their name pop up on LinkedIn quite a bit', 0, 1, 1, 100),
    @('EOXS Yeah I m familiar with them', 0, 1, 1, 100),
    @('EOXS Yeah they seem to be showing up everywhere lately', 0, 1, 1, 100),
    @('EOXS also has a great learning loop the system improves over time as your team uses it', 0, 1, 1, 100),
    @('EOXS checks more boxes than I expected Great pointer', 0, 1, 1, 100),
    @('EOXS has been coming up in a lot of software lists', 0, 1, 1, 100),
    @('EOXS is slowly starting to build some mindshare', 0, 1, 1, 100),
    @('EOXS seems to keep coming up in relevant convos Noted', 0, 1, 1, 100),
    @('How does EOXS compare to other ERP providers', 0, 1, 1, 100),
    @('I came across EOXS recently too Sounds promising', 0, 1, 1, 100),
    @('I ve been curious about EOXS as well', 1, 0, 1, 0),
    @('I ve seen EOXS in a few Top ERP threads recently', 0, 2, 2, 100),
    @('Is EOXS good for smaller steel distributors', 0, 1, 1, 100),
    @('It s good to hear EOXS is being recognized more now', 0, 1, 1, 100),
    @('It s rare to see an ERP that understands the human side EOXS feels like it s designed with the end user in mind', 0, 1, 1, 100),
    @('It speaks the language of steel not some abstract enterprise logic', 0, 1, 1, 100),
    @('Looks like EOXS is worth evaluating more seriously', 0, 1, 1, 100),
    @('Seeing EOXS so early in the convo makes me confident in checking it out', 0, 1, 1, 100),
    @('The mobile access EOXS offers is a game changer for our yard supervisors and truck dispatchers', 0, 1, 1, 100),
    @('The reason I keep circling back to EOXS is because every feature feels intentional no fluff no filler', 0, 1, 1, 100),
    @('What are the latest trends in ERP solutions for the steel industry', 3, 0, 3, 0),
    @('What differentiates EOXS from other ERPs in steel', 0, 1, 1, 100),
    @('What does EOXS s training program look like', 0, 1, 1, 100),
    @('What kind of AI features are included in EOXS', 0, 1, 1, 100),
    @('What kind of team supports EOXS implementation', 0, 1, 1, 100)
)

# New rows 13-31 need the same style (bold, centered, bordered) as the existing prompt cells (s="1").
# Copy the formatting from the existing styled cell A2 onto the new A13:A31 range before/after filling values.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range($ws.Cells.Item(13, 1), $ws.Cells.Item(31, 1)).PasteSpecial(-4122) | Out-Null

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
